$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Goias vs Vila Nova FC
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 4.75
$ws.Range("J3").Value = 2.5
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 5.5
$ws.Range("N3").Value = 7.5
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 7.5
$ws.Range("Z3").Value = 13
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 17
$ws.Range("AN3").Value = 3.6
$ws.Range("AQ3").Value = 34
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 29
$ws.Range("BA3").Value = 151

# Row 5: Aguilas vs Dep. Cali
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 2.75
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("U5").Value = 2.38
$ws.Range("V5").Value = 1.53
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 7.5
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 17
$ws.Range("AJ5").Value = 51
$ws.Range("AU5").Value = 10
$ws.Range("AW5").Value = 6

# Row 7: Leones Negros vs Atletico La Paz
$ws.Range("G7").Value = 1.39
$ws.Range("H7").Value = 4.5
$ws.Range("I7").Value = 6.8
$ws.Range("J7").Value = 1.85
$ws.Range("K7").Value = 2.4
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 12.9
$ws.Range("O7").Value = 1.19
$ws.Range("P7").Value = 3.65
$ws.Range("Q7").Value = 1.6
$ws.Range("R7").Value = 2.07
$ws.Range("S7").Value = 1.3
$ws.Range("T7").Value = 3.32
$ws.Range("U7").Value = 1.82
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 7.5
$ws.Range("X7").Value = 6.9
$ws.Range("Z7").Value = 9.25
$ws.Range("AA7").Value = 11
$ws.Range("AB7").Value = 25
$ws.Range("AC7").Value = 13
$ws.Range("AD7").Value = 9
$ws.Range("AE7").Value = 19
$ws.Range("AG7").Value = 19
$ws.Range("AH7").Value = 45
$ws.Range("AL7").Value = 65
$ws.Range("AO7").Value = 6.2
$ws.Range("AP7").Value = 15.5
$ws.Range("AQ7").Value = 17
$ws.Range("AS7").Value = 200
$ws.Range("AT7").Value = 3.1
$ws.Range("AW7").Value = 8
$ws.Range("AY7").Value = 40
$ws.Range("BB7").Value = 500

# Row 8: Alebrijes Oaxaca vs Celaya
$ws.Range("G8").Value = 3.65
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 1.91
$ws.Range("J8").Value = 4.1
$ws.Range("K8").Value = 2.1
$ws.Range("L8").Value = 2.5
$ws.Range("N8").Value = 9.85
$ws.Range("P8").Value = 3.15
$ws.Range("Q8").Value = 1.82
$ws.Range("R8").Value = 1.88
$ws.Range("S8").Value = 1.37
$ws.Range("T8").Value = 2.94
$ws.Range("U8").Value = 1.65
$ws.Range("V8").Value = 1.98
$ws.Range("W8").Value = 11
$ws.Range("X8").Value = 20
$ws.Range("Y8").Value = 12.5
$ws.Range("AA8").Value = 32
$ws.Range("AB8").Value = 37
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 6.7
$ws.Range("AG8").Value = 7.8
$ws.Range("AH8").Value = 9.5
$ws.Range("AJ8").Value = 17
$ws.Range("AK8").Value = 14.5
$ws.Range("AN8").Value = 5.5
$ws.Range("AS8").Value = 300
$ws.Range("AT8").Value = 2.7
$ws.Range("AV8").Value = 65
$ws.Range("AW8").Value = 3.8
$ws.Range("AX8").Value = 9.5
$ws.Range("AY8").Value = 18
$ws.Range("AZ8").Value = 35
$ws.Range("BA8").Value = 65

# Row 11: Cerro Porteno vs General Caballero JLM
$ws.Range("G11").Value = 1.38
$ws.Range("H11").Value = 4.33
$ws.Range("I11").Value = 7
$ws.Range("J11").Value = 1.95
$ws.Range("K11").Value = 2.25
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 9.5
$ws.Range("U11").Value = 2.25
$ws.Range("V11").Value = 1.57
$ws.Range("X11").Value = 6
$ws.Range("AD11").Value = 9
$ws.Range("AF11").Value = 81
$ws.Range("AO11").Value = 7
$ws.Range("AW11").Value = 8.5
$ws.Range("BA11").Value = 201

# Row 12: Academico Viseu vs Maritimo
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 1.36
$ws.Range("P12").Value = 3
$ws.Range("Q12").Value = 2.15
$ws.Range("R12").Value = 1.67

# Row 13: Pacos Ferreira vs Torreense
$ws.Range("G13").Value = 2.15
$ws.Range("I13").Value = 3.3
$ws.Range("J13").Value = 2.88
$ws.Range("K13").Value = 2.05
$ws.Range("O13").Value = 1.36
$ws.Range("P13").Value = 3
$ws.Range("Q13").Value = 2.15
$ws.Range("R13").Value = 1.67
$ws.Range("X13").Value = 10
$ws.Range("Y13").Value = 9.5
$ws.Range("Z13").Value = 21
$ws.Range("AG13").Value = 9.5
$ws.Range("AI13").Value = 12
$ws.Range("AJ13").Value = 34
$ws.Range("AS13").Value = 201
$ws.Range("AW13").Value = 5
$ws.Range("AX13").Value = 19
$ws.Range("BA13").Value = 81

# Row 14: Chaves vs Feirense
$ws.Range("G14").Value = 1.85
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 3.9
$ws.Range("J14").Value = 2.6
$ws.Range("L14").Value = 4.5
$ws.Range("X14").Value = 8.5
$ws.Range("AA14").Value = 17
$ws.Range("AI14").Value = 13
$ws.Range("AK14").Value = 34

# Row 15: Huesca vs Albacete
$ws.Range("G15").Value = 2.25
$ws.Range("I15").Value = 3.4
$ws.Range("J15").Value = 3
$ws.Range("W15").Value = 6.5
$ws.Range("Y15").Value = 9.5
$ws.Range("AG15").Value = 9

# Row 16: R. Oviedo vs Almeria
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 9
$ws.Range("Q16").Value = 2.1
$ws.Range("R16").Value = 1.7
$ws.Range("AS16").Value = 151
$ws.Range("AX16").Value = 19

# Row 18: Burgos CF vs Mirandes
$ws.Range("H18").Value = 2.8
$ws.Range("I18").Value = 3.75
$ws.Range("J18").Value = 3.2
$ws.Range("K18").Value = 1.83
$ws.Range("L18").Value = 4.5
$ws.Range("M18").Value = 1.14
$ws.Range("N18").Value = 5.5
$ws.Range("O18").Value = 1.67
$ws.Range("P18").Value = 2.1
$ws.Range("Q18").Value = 3.1
$ws.Range("R18").Value = 1.36
$ws.Range("S18").Value = 1.67
$ws.Range("T18").Value = 2.1
$ws.Range("U18").Value = 2.38
$ws.Range("V18").Value = 1.53
$ws.Range("Y18").Value = 11
$ws.Range("AA18").Value = 26
$ws.Range("AF18").Value = 101
$ws.Range("AG18").Value = 7.5
$ws.Range("AO18").Value = 15
$ws.Range("AP18").Value = 34
$ws.Range("AR18").Value = 101
$ws.Range("AS18").Value = 351
$ws.Range("AT18").Value = 2.1
$ws.Range("AU18").Value = 10
$ws.Range("AX18").Value = 23
$ws.Range("BA18").Value = 151

# Row 20: Elche vs Dep. La Coruna
$ws.Range("G20").Value = 1.95
$ws.Range("H20").Value = 3.1
$ws.Range("I20").Value = 4.33
$ws.Range("J20").Value = 2.75
$ws.Range("L20").Value = 4.5
$ws.Range("Q20").Value = 2.2
$ws.Range("R20").Value = 1.65
$ws.Range("S20").Value = 1.44
$ws.Range("T20").Value = 2.63
$ws.Range("U20").Value = 1.91
$ws.Range("V20").Value = 1.8
$ws.Range("W20").Value = 6.5
$ws.Range("X20").Value = 8.5
$ws.Range("Y20").Value = 9
$ws.Range("Z20").Value = 17
$ws.Range("AA20").Value = 17
$ws.Range("AD20").Value = 6
$ws.Range("AG20").Value = 11
$ws.Range("AH20").Value = 21
$ws.Range("AI20").Value = 15
$ws.Range("AK20").Value = 41
$ws.Range("AM20").Value = 351
$ws.Range("AO20").Value = 11
$ws.Range("AT20").Value = 2.63
$ws.Range("AX20").Value = 23
$ws.Range("AY20").Value = 34
$ws.Range("AZ20").Value = 81

# Row 22: Cerro Largo vs Nacional
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 3.6
$ws.Range("I22").Value = 1.53
$ws.Range("M22").Value = 1.08
$ws.Range("N22").Value = 8
$ws.Range("W22").Value = 15
$ws.Range("Z22").Value = 81
$ws.Range("AF22").Value = 67
$ws.Range("AR22").Value = 201
$ws.Range("AY22").Value = 21

# Row 27: Universidad Central vs Caracas
$ws.Range("H27").Value = 2.9
$ws.Range("I27").Value = 3.5
$ws.Range("K27").Value = 1.91
$ws.Range("L27").Value = 4.1
$ws.Range("N27").Value = 6.3
$ws.Range("P27").Value = 2.4
$ws.Range("Q27").Value = 2.27
$ws.Range("S27").Value = 1.5
$ws.Range("T27").Value = 2.27
$ws.Range("AC27").Value = 6.7
$ws.Range("AD27").Value = 5.7
$ws.Range("AG27").Value = 8
$ws.Range("AH27").Value = 17
$ws.Range("AI27").Value = 12.5
$ws.Range("AK27").Value = 37
$ws.Range("AR27").Value = 90
$ws.Range("AT27").Value = 2.25
$ws.Range("AU27").Value = 7.4
$ws.Range("AW27").Value = 5.2
$ws.Range("AX27").Value = 20
$ws.Range("AY27").Value = 30
$ws.Range("BA27").Value = 175
$ws.Range("BB27").Value = 450
